$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 10.439420823763685
$ws.Range("C2").Value = -0.27451878770839427
$ws.Range("D2").Value = 0.37687517840855428
$ws.Range("E2").Value = 0.1343756834943152

# Row 3 values
$ws.Range("B3").Value = 3.2310767151311972
$ws.Range("C3").Value = 5.1597658033607727
$ws.Range("D3").Value = 2.3353231463272066
$ws.Range("E3").Value = -0.59635159644502955

# Update selection to match new active range B1:E3
$ws.Range("B1:E3").Select()
